$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A22").Value = "sekundy najedzonego to 570454"
$ws.Range("A22").Select()
